$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.356.02'
$ws.Cells.Item(2, 5).Value = '  +2.31%  '

$ws.Cells.Item(3, 4).Value = '1.662.15'
$ws.Cells.Item(3, 5).Value = '  +1.30%  '

$ws.Cells.Item(4, 5).Value = '  -0.49%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '220.15'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.11%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.507'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.80%  '

$ws.Cells.Item(7, 5).Value = '  -0.46%  '

$ws.Cells.Item(8, 5).Value = '  +1.26%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.0627'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.33%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '20.04'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +4.80%  '

$ws.Cells.Item(11, 5).Value = '  +0.82%  '

$ws.Cells.Item(12, 4).Value = '1.894.66'
$ws.Cells.Item(12, 5).Value = '  +1.31%  '

$ws.Cells.Item(13, 4).Value = '1.660.36'
$ws.Cells.Item(13, 5).Value = '  +1.14%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.20'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +0.96%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.533'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.22%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '67.21'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +3.93%  '

$ws.Cells.Item(17, 4).Value = '27.342.21'
$ws.Cells.Item(17, 5).Value = '  +2.23%  '

$ws.Cells.Item(18, 4).Value = '0.0₃0737'
$ws.Cells.Item(18, 5).Value = '  +0.55%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '222.77'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +3.99%  '

$ws.Cells.Item(20, 5).Value = '  -0.47%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.77'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +8.90%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.45'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.70%  '

$ws.Cells.Item(23, 5).Value = '  +5.55%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '9.29'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +0.22%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '147.02'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.94%  '

$ws.Cells.Item(26, 5).Value = '  -0.45%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.43'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +3.80%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.119'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +0.96%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '16.07'
$ws.Cells.Item(29, 4).Style = 'Normal'

$ws.Cells.Item(30, 5).Value = '  +1.20%  '

$ws.Cells.Item(31, 5).Value = '  +0.93%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.39'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +0.15%  '

$ws.Cells.Item(33, 5).Value = '  +0.01%  '

$ws.Cells.Item(34, 5).Value = '  +2.36%  '

$ws.Cells.Item(35, 4).Value = '1.266.42'
$ws.Cells.Item(35, 5).Value = '  -1.61%  '

$ws.Cells.Item(36, 5).Value = '  +0.74%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.0178'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.39%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.538'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +0.23%  '

$ws.Cells.Item(39, 5).Value = '  +2.37%  '

$ws.Cells.Item(40, 5).Value = '  -0.43%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.813'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +0.88%  '

$ws.Cells.Item(42, 5).Value = '  +2.52%  '

$ws.Cells.Item(43, 4).Value = '1.805.69'
$ws.Cells.Item(43, 5).Value = '  +1.46%  '

$ws.Cells.Item(44, 5).Value = '  -4.03%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '61.77'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +1.29%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '92.38'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.81%  '

$ws.Cells.Item(47, 5).Value = '  +1.74%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0518'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.16%  '

$ws.Cells.Item(49, 5).Value = '  +2.07%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '7.70'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.36%  '
